$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 15/16: Chainlink and WrappedEther swap places (name + link together)
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

# Updated Price (D) and Volume(1h) (E) values for each coin row.
# The source data stores these as plain text, not numbers, so for any
# value that Excel would otherwise auto-convert to a number (e.g. "1.001"),
# a leading apostrophe forces literal text; re-applying the Normal style
# afterwards clears the "stored as text" quote-prefix formatting flag that
# the apostrophe trick leaves behind, so cell formatting is left untouched.
$ws.Range("D2").Value = "28.311.43"
$ws.Range("E2").Value = "  +5.23%  "
$ws.Range("D3").Value = "1.802.84"
$ws.Range("E3").Value = "  +3.76%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'316.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.5483"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.12%  "
$ws.Range("D8").Value = "'0.3841"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.75%  "
$ws.Range("D9").Value = "'42.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").Value = "'0.07566"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.41%  "
$ws.Range("E11").Value = "  +5.91%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "'21.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.67%  "
$ws.Range("D14").Value = "'6.199"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.40%  "
$ws.Range("D15").Value = "'7.332"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.61%  "
$ws.Range("D16").Value = "1.797.16"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "'91.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.73%  "
$ws.Range("D18").Value = "'0.00001069"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.45%  "
$ws.Range("D19").Value = "'0.06450"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'17.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.70%  "
$ws.Range("D22").Value = "'5.975"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.64%  "
$ws.Range("D23").Value = "28.307.24"
$ws.Range("E23").Value = "  +4.96%  "
$ws.Range("D24").Value = "'11.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "'2.174"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.10%  "
$ws.Range("D26").Value = "'158.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").Value = "'20.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.99%  "
$ws.Range("D28").Value = "'2.428"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.98%  "
$ws.Range("D29").Value = "2.008.56"
$ws.Range("E29").Value = "  +3.53%  "
$ws.Range("D30").Value = "'123.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.98%  "
$ws.Range("D31").Value = "'1.144"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.66%  "
$ws.Range("D32").Value = "'0.1016"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.86%  "
$ws.Range("D33").Value = "'5.728"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.02%  "
$ws.Range("D34").Value = "'3.648"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").Value = "'0.2312"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.58%  "
$ws.Range("D36").Value = "'0.06304"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.40%  "
$ws.Range("D37").Value = "'8.921"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +17.73%  "
$ws.Range("D38").Value = "'0.02319"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.58%  "
$ws.Range("D39").Value = "'11.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.74%  "
$ws.Range("D40").Value = "'5.023"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.28%  "
$ws.Range("D41").Value = "'0.6371"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.54%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'1.157"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.22%  "
$ws.Range("D44").Value = "'1.384"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").Value = "'13.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.79%  "
$ws.Range("D46").Value = "'0.5978"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.86%  "
$ws.Range("D47").Value = "'3.678"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("D48").Value = "'124.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").Value = "'1.972"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.52%  "
$ws.Range("D50").Value = "'1.146"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("D51").Value = "'0.06892"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.37%  "
